# Updates cryptos list data (prices and volume percentages) as of the
# commit "Updated cryptos list on Mon Jan 29 05:41:59 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.166.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.265.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.527'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.11%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.617.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.266.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.055.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.49%  '
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.57%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.13%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0737'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("E38").Value = '  -4.33%  '
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.46%  '
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.949.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '92.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("E51").Value = '  -2.57%  '
